$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cluster_12")

# --- Restructure columns ---
# Before: A=result, B=processed_result, C=category
# After:  A=id, B=category, C=severity, D=recurrent, E=result, F=processed_result

# Insert 2 blank columns before column A.
# Layout becomes: A=(blank), B=(blank), C=result, D=processed_result, E=category
$ws.Range("A1:B1").EntireColumn.Insert()

# Move category (E) into B, freeing E.
$ws.Range("E1:E32").Cut($ws.Range("B1"))

# Move processed_result (D) into F, freeing D.
$ws.Range("D1:D32").Cut($ws.Range("F1"))

# Move result (C) into E, freeing C.
$ws.Range("C1:C32").Cut($ws.Range("E1"))

# Now: A=(blank,id target), B=category, C=(blank,severity target), D=(blank,recurrent target), E=result, F=processed_result

# Copy the bold/centered header style from B1 onto A1, C1 and D1 (the newly introduced header columns).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Set header text ---
$ws.Range("A1").Value = "id"
$ws.Range("C1").Value = "severity"
$ws.Range("D1").Value = "recurrent"

# --- Populate id / severity / recurrent for each data row ---
$ws.Range("A2").Value = 10010000007514
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 5
$ws.Range("A3").Value = 10010000007615
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 5
$ws.Range("A4").Value = 10010000007651
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5
$ws.Range("A5").Value = 10010000007643
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 5
$ws.Range("A6").Value = 10010000007892
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 5
$ws.Range("A7").Value = 10010000008834
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 5
$ws.Range("A8").Value = 10010000007895
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("A9").Value = 10010000007951
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 5
$ws.Range("A10").Value = 10010000008326
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 4
$ws.Range("A11").Value = 10010000007984
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 5
$ws.Range("A12").Value = 10010000008005
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 5
$ws.Range("A13").Value = 10010000008404
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 5
$ws.Range("A14").Value = 10010000008338
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 5
$ws.Range("A15").Value = 10010000008380
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 5
$ws.Range("A16").Value = 10010000008482
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("A17").Value = 10010000008596
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("A18").Value = 10010000008227
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("A19").Value = 10010000008851
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 5
$ws.Range("A20").Value = 10010000008719
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("A21").Value = 10010000009472
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 5
$ws.Range("A22").Value = 10010000009727
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 5
$ws.Range("A23").Value = 10010000009919
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("A24").Value = 10010000009323
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("A25").Value = 10010000009357
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 5
$ws.Range("A26").Value = 10010000007803
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("A27").Value = 10010000007965
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 5
$ws.Range("A28").Value = 10010000007994
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 5
$ws.Range("A29").Value = 10010000008277
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 5
$ws.Range("A30").Value = 10010000009617
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 5
$ws.Range("A31").Value = 10010000009466
$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 5
$ws.Range("A32").Value = 10010000009916
$ws.Range("C32").Value = 4
$ws.Range("D32").Value = 4

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
